# Auto-generated Excel COM-interop script
# Applies scraped market price / profit updates across multiple sheets
# (ALC, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled runner diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1250600
$ws.Range("J2").Value = 799.8
$ws.Range("L2").Value = 799.8
$ws.Range("N2").Value = -1025.8
$ws.Range("H70").Value = 2030
$ws.Range("I70").Value = 995
$ws.Range("J70").Value = 4100
$ws.Range("K70").Value = 2985
$ws.Range("L70").Value = 12300
$ws.Range("M70").Value = -2715
$ws.Range("N70").Value = -12840
$ws.Range("H73").Value = 2030
$ws.Range("I73").Value = 995
$ws.Range("J73").Value = 4100
$ws.Range("K73").Value = 2985
$ws.Range("L73").Value = 12300
$ws.Range("M73").Value = -2049
$ws.Range("N73").Value = -14172
$ws.Range("H129").Value = 1053.7142
$ws.Range("I129").Value = 369
$ws.Range("J129").Value = 1167.8334
$ws.Range("K129").Value = 1107
$ws.Range("L129").Value = 3503.5002
$ws.Range("M129").Value = 3893
$ws.Range("N129").Value = -13503.5002
$ws.Range("H132").Value = 239622.45
$ws.Range("I132").Value = 1493.1351
$ws.Range("J132").Value = 2001779.4
$ws.Range("K132").Value = 4479.4053
$ws.Range("L132").Value = 6005338.199999999
$ws.Range("M132").Value = -1949.4053
$ws.Range("N132").Value = -6010398.199999999
$ws.Range("H138").Value = 3654.79
$ws.Range("I138").Value = 814.7
$ws.Range("J138").Value = 4364.8125
$ws.Range("K138").Value = 2444.1
$ws.Range("L138").Value = 13094.4375
$ws.Range("M138").Value = 2695.9
$ws.Range("N138").Value = -23374.4375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 56796.188
$ws.Range("I134").Value = 3059.6785
$ws.Range("J134").Value = 223976.44
$ws.Range("K134").Value = 9179.0355
$ws.Range("L134").Value = 671929.3200000001
$ws.Range("M134").Value = -6644.0355
$ws.Range("N134").Value = -676999.3200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6332.9546
$ws.Range("I58").Value = 2193.7273
$ws.Range("J58").Value = 10472.182
$ws.Range("K58").Value = 2193.7273
$ws.Range("L58").Value = 10472.182
$ws.Range("M58").Value = -1990.7273
$ws.Range("N58").Value = -10878.182
$ws.Range("H94").Value = 3330.25
$ws.Range("I94").Value = 3356
$ws.Range("J94").Value = 3321.6667
$ws.Range("K94").Value = 3356
$ws.Range("L94").Value = 3321.6667
$ws.Range("M94").Value = -2905
$ws.Range("N94").Value = -4223.6667
$ws.Range("H136").Value = 6332.9546
$ws.Range("I136").Value = 2193.7273
$ws.Range("J136").Value = 10472.182
$ws.Range("K136").Value = 6581.1819
$ws.Range("L136").Value = 31416.546
$ws.Range("M136").Value = -4031.1819
$ws.Range("N136").Value = -36516.546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2907.5
$ws.Range("I56").Value = 2907.5
$ws.Range("K56").Value = 2907.5
$ws.Range("M56").Value = -2377.5
$ws.Range("H113").Value = 518.3333
$ws.Range("J113").Value = 526.6667
$ws.Range("L113").Value = 1580.0001
$ws.Range("N113").Value = -5920.0001
$ws.Range("H117").Value = 2700.1538
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 2841.8333
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 8525.499899999999
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -15409.4999
$ws.Range("H129").Value = 2346.6333
$ws.Range("I129").Value = 1274.5454
$ws.Range("J129").Value = 2967.3157
$ws.Range("K129").Value = 3823.6362
$ws.Range("L129").Value = 8901.947100000001
$ws.Range("M129").Value = 1176.3638
$ws.Range("N129").Value = -18901.9471
$ws.Range("H136").Value = 3113.75
$ws.Range("I136").Value = 1992
$ws.Range("J136").Value = 4983.3335
$ws.Range("K136").Value = 5976
$ws.Range("L136").Value = 14950.0005
$ws.Range("M136").Value = -876
$ws.Range("N136").Value = -25150.0005
$ws.Range("H138").Value = 2942.875
$ws.Range("I138").Value = 2585
$ws.Range("K138").Value = 7755
$ws.Range("M138").Value = -2615

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1373.6
$ws.Range("I113").Value = 1258.5714
$ws.Range("J113").Value = 1642
$ws.Range("K113").Value = 1258.5714
$ws.Range("L113").Value = 1642
$ws.Range("M113").Value = 911.4286
$ws.Range("N113").Value = -5982
$ws.Range("H122").Value = 2283.3872
$ws.Range("I122").Value = 1763.6471
$ws.Range("J122").Value = 2914.5
$ws.Range("K122").Value = 5290.9413
$ws.Range("L122").Value = 8743.5
$ws.Range("M122").Value = -2840.9413
$ws.Range("N122").Value = -13643.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1090.8334
$ws.Range("J22").Value = 1314.2858
$ws.Range("L22").Value = 1314.2858
$ws.Range("N22").Value = -1904.2858
$ws.Range("H27").Value = 1090.8334
$ws.Range("J27").Value = 1314.2858
$ws.Range("L27").Value = 1314.2858
$ws.Range("N27").Value = -1528.2858
$ws.Range("H46").Value = 1300
$ws.Range("I46").Value = 2300
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 2300
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -2112
$ws.Range("N46").Value = -1276
$ws.Range("H68").Value = 3950.25
$ws.Range("I68").Value = 4200.4
$ws.Range("J68").Value = 3533.3333
$ws.Range("K68").Value = 4200.4
$ws.Range("L68").Value = 3533.3333
$ws.Range("M68").Value = -3451.4
$ws.Range("N68").Value = -5031.3333
$ws.Range("H71").Value = 3950.25
$ws.Range("I71").Value = 4200.4
$ws.Range("J71").Value = 3533.3333
$ws.Range("K71").Value = 21002
$ws.Range("L71").Value = 17666.6665
$ws.Range("M71").Value = -17258
$ws.Range("N71").Value = -25154.6665
$ws.Range("H94").Value = 20000
$ws.Range("I94").Value = 20000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 20000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -19324
$ws.Range("N94").ClearContents()
$ws.Range("H136").Value = 2758.1155
$ws.Range("I136").Value = 1234.8096
$ws.Range("J136").Value = 9156
$ws.Range("K136").Value = 3704.4288
$ws.Range("L136").Value = 27468
$ws.Range("M136").Value = -1154.4288
$ws.Range("N136").Value = -32568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1718.1852
$ws.Range("I136").Value = 1695.64
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5086.92
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2536.92
$ws.Range("N136").Value = -11100
